$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B11").Value = "cf2f1f7db0e2df3878213475336921a3"
$ws.Range("B15").Value = "8add17134bbad18b7acdadc0e8ee3441"
$ws.Range("B29").Value = "020895d6d25807ceea798e4b8cd227b7"
$ws.Range("B121").Value = "4738ec7c82e727e7534c7a7fb08cb771"
$ws.Range("B126").Value = "cb2b48530b102a7818d954df99d33a88"
$ws.Range("B133").Value = "a1f0aeac8802c0250624fc9fa2c26529"
$ws.Range("B159").Value = "567cf77756c9ad1d2efe5d8d378938af"
$ws.Range("B162").Value = "f6e10bcb8d47e08bc2d03119866ebb46"
$ws.Range("B169").Value = "bd9aecd057b8b0de503941ed6157bbd2"
$ws.Range("B175").Value = "2ac35ee76222d5df5e0654457b0c544f"
$ws.Range("B191").Value = "32cfcd119d179ac4a5597dc259240032"
$ws.Range("B198").Value = "218495e0ce5c193e6de4326bb103aa11"
$ws.Range("B293").Value = "d9e41eccb1727d9b81e0c2d1587a1e06"
$ws.Range("B302").Value = "2adb940599b723985c03c239ef449b0e"
$ws.Range("B339").Value = "95dbda5d9a8b6ad8dfae2c4599d555fd"
$ws.Range("B420").Value = "bf3569543f5afe0bd329968445d710df"
$ws.Range("B485").Value = "137211fbc02800389c315e3667e2f3e3"
$ws.Range("B506").Value = "c4e086901e87a390d81c08e4bb9fdebd"
$ws.Range("B508").Value = "e5689301a7dcef202aae3ff556c77d8e"
$ws.Range("B558").Value = "98a5c43dfa1645b5e2a64db62b84bf4f"
$ws.Range("B624").Value = "49317de9592d0ba2745f2811467e0469"
$ws.Range("B635").Value = "ba40ada3c09dc5fca60bdcd71f8eb628"
$ws.Range("B637").Value = "93ef2328a3b5c2a9f75453d8c4ad9cbd"
$ws.Range("B657").Value = "13b07137e0f076a52f2d182dfef486cd"
$ws.Range("B663").Value = "7aa8a8d7a5e511b7f5357c779c3135d0"
$ws.Range("B673").Value = "09702f670aedd0e8bb30def8cf4aec9c"
$ws.Range("B688").Value = "02796346b86ff6d9d6c7fce4bac0cac5"
$ws.Range("B708").Value = "ac11ff5172c43564a5b15233fd7c3275"
$ws.Range("B711").Value = "2bbbc64dc8be0d94d0befb3fe111fabd"
$ws.Range("B741").Value = "c406e93abe460dfbf507cba21d7187c5"
$ws.Range("B827").Value = "8984ed957ef45588ab2b7e250414079d"
$ws.Range("B870").Value = "73ac72d57a94466bf0648eef63be2fea"
